$wb = $excel.ActiveWorkbook

# Sheet references (by index, order matches workbook.xml: 展览, 演出, 本地生活, 全部类型)
$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibitions)
$ws2 = $wb.Worksheets.Item(2)   # 演出 (Performances)
$ws3 = $wb.Worksheets.Item(3)   # 本地生活 (Local Life)
$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All Types - aggregated)

$ws1.Range("F2").Value = 48
$ws1.Range("F6").Value = 1698
$ws1.Range("F9").Value = 2470
$ws1.Range("F10").Value = 714
$ws1.Range("F13").Value = 15
$ws1.Range("F16").Value = 209
$ws1.Range("F20").Value = 701
$ws1.Range("F24").Value = 24
$ws1.Range("D27").Value = "淮海中路775号 niko and ......"
$ws1.Range("F27").Value = 32
$ws1.Range("F29").Value = 1760
$ws1.Range("F31").Value = 530
$ws1.Range("F32").Value = 518
$ws1.Range("F34").Value = 81
$ws1.Range("F35").Value = 4546
$ws1.Range("F36").Value = 119
$ws2.Range("F3").Value = 384
$ws2.Range("F11").Value = 64
$ws2.Range("F15").Value = 318
$ws2.Range("F21").Value = 2
$ws2.Range("F25").Value = 1761
$ws2.Range("F26").Value = 242
$ws2.Range("F37").Value = 65
$ws3.Range("F4").Value = 1413
$ws3.Range("F5").Value = 1784
$ws3.Range("F7").Value = 177
$ws4.Range("F3").Value = 1413
$ws4.Range("F4").Value = 384
$ws4.Range("F6").Value = 48
$ws4.Range("F10").Value = 1698
$ws4.Range("G12").Value = 90
$ws4.Range("F16").Value = 2470
$ws4.Range("F17").Value = 714
$ws4.Range("F20").Value = 15
$ws4.Range("F23").Value = 64
$ws4.Range("F25").Value = 209
$ws4.Range("F26").Value = 318
$ws4.Range("F36").Value = 177
$ws4.Range("F38").Value = 1761
$ws4.Range("F40").Value = 1760
$ws4.Range("F41").Value = 242
$ws4.Range("F42").Value = 518
$ws4.Range("F45").Value = 4546
$ws4.Range("F46").Value = 119
